$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 33.06592966666667
$ws.Range("H2").Value = 99.197789
$ws.Range("I2").Value = 0.4620579289161133
$ws.Range("J2").Value = 0.4620579289161132
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.023694333333334
$ws.Range("N2").Value = 21.071083
$ws.Range("O2").Value = 0.1590811435055747
$ws.Range("P2").Value = 0.1590811435055747
$ws.Range("Q2").Value = 232.2449828261653
$ws.Range("R2").Value = 2090.204845435487
$ws.Range("S2").Value = 0.07350470369779287
$ws.Range("T2").Value = 0.07350470369779283
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 33.06592966666667
$ws.Range("H3").Value = 99.197789
$ws.Range("I3").Value = 0.4620579289161133
$ws.Range("J3").Value = 0.4620579289161132
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.91445766666667
$ws.Range("N3").Value = 41.74337300000001
$ws.Range("O3").Value = 0.3151515045818828
$ws.Range("P3").Value = 0.3151515045818827
$ws.Range("Q3").Value = 460.0944785558109
$ws.Range("R3").Value = 4140.850307002298
$ws.Range("S3").Value = 0.1456182515019018
$ws.Range("T3").Value = 0.1456182515019017
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 33.06592966666667
$ws.Range("H4").Value = 99.197789
$ws.Range("I4").Value = 0.4620579289161133
$ws.Range("J4").Value = 0.4620579289161132
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.213494
$ws.Range("N4").Value = 69.640482
$ws.Range("O4").Value = 0.5257673519125425
$ws.Range("P4").Value = 0.5257673519125424
$ws.Range("Q4").Value = 767.5757599215888
$ws.Range("R4").Value = 6908.181839294299
$ws.Range("S4").Value = 0.2429349737164187
$ws.Range("T4").Value = 0.2429349737164186
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.69090766666667
$ws.Range("H5").Value = 44.072723
$ws.Range("I5").Value = 0.2052883568914378
$ws.Range("J5").Value = 0.2052883568914378
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.023694333333334
$ws.Range("N5").Value = 21.071083
$ws.Range("O5").Value = 0.1590811435055747
$ws.Range("P5").Value = 0.1590811435055747
$ws.Range("Q5").Value = 103.1844449298899
$ws.Range("R5").Value = 928.660004369009
$ws.Range("S5").Value = 0.03265750656267046
$ws.Range("T5").Value = 0.03265750656267045
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.69090766666667
$ws.Range("H6").Value = 44.072723
$ws.Range("I6").Value = 0.2052883568914378
$ws.Range("J6").Value = 0.2052883568914378
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.91445766666667
$ws.Range("N6").Value = 41.74337300000001
$ws.Range("O6").Value = 0.3151515045818828
$ws.Range("P6").Value = 0.3151515045818827
$ws.Range("Q6").Value = 204.4160128127421
$ws.Range("R6").Value = 1839.744115314679
$ws.Range("S6").Value = 0.06469693454747917
$ws.Range("T6").Value = 0.06469693454747916
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.69090766666667
$ws.Range("H7").Value = 44.072723
$ws.Range("I7").Value = 0.2052883568914378
$ws.Range("J7").Value = 0.2052883568914378
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.213494
$ws.Range("N7").Value = 69.640482
$ws.Range("O7").Value = 0.5257673519125425
$ws.Range("P7").Value = 0.5257673519125424
$ws.Range("Q7").Value = 341.0272969747207
$ws.Range("R7").Value = 3069.245672772486
$ws.Range("S7").Value = 0.1079339157812882
$ws.Range("T7").Value = 0.1079339157812882
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 23.805466
$ws.Range("H8").Value = 71.416398
$ws.Range("I8").Value = 0.3326537141924489
$ws.Range("J8").Value = 0.3326537141924489
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.023694333333334
$ws.Range("N8").Value = 21.071083
$ws.Range("O8").Value = 0.1590811435055747
$ws.Range("P8").Value = 0.1590811435055747
$ws.Range("Q8").Value = 167.2023166465594
$ws.Range("R8").Value = 1504.820849819034
$ws.Range("S8").Value = 0.0529189332451114
$ws.Range("T8").Value = 0.05291893324511138
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 23.805466
$ws.Range("H9").Value = 71.416398
$ws.Range("I9").Value = 0.3326537141924489
$ws.Range("J9").Value = 0.3326537141924489
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.91445766666667
$ws.Range("N9").Value = 41.74337300000001
$ws.Range("O9").Value = 0.3151515045818828
$ws.Range("P9").Value = 0.3151515045818827
$ws.Range("Q9").Value = 331.2401488922727
$ws.Range("R9").Value = 2981.161340030455
$ws.Range("S9").Value = 0.1048363185325019
$ws.Range("T9").Value = 0.1048363185325019
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 23.805466
$ws.Range("H10").Value = 71.416398
$ws.Range("I10").Value = 0.3326537141924489
$ws.Range("J10").Value = 0.3326537141924489
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.213494
$ws.Range("N10").Value = 69.640482
$ws.Range("O10").Value = 0.5257673519125425
$ws.Range("P10").Value = 0.5257673519125424
$ws.Range("Q10").Value = 552.608042158204
$ws.Range("R10").Value = 4973.472379423836
$ws.Range("S10").Value = 0.1748984624148356
$ws.Range("T10").Value = 0.1748984624148356
